# Week 13 logging update — Target Depth Data (Seahawks)
# Updates the "H" (home) row totals on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 2 ("H") ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 374
$wsOff.Range("C2").Value = 274
$wsOff.Range("D2").Value = 93
$wsOff.Range("E2").Value = 43
$wsOff.Range("F2").Value = 5

# --- DEF sheet: row 2 ("H") ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 567
$wsDef.Range("C2").Value = 416
$wsDef.Range("D2").Value = 109
$wsDef.Range("E2").Value = 51
$wsDef.Range("F2").Value = 3
$wsDef.Range("G2").Value = 7
